# Rename the "ATP-Fs" category to "ATP-Finals" (big-data-import cleanup).
# This is the block of rows describing the ATP Finals (PlayersNumber = 8)
# category in the Sheet1 points table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 46; $row -le 56; $row++) {
    $ws.Cells.Item($row, 1).Value = "ATP-Finals"
}

# Leave the view scrolled/selected where the edit happened.
$ws.Range("A40").Select()
try {
    $excel.ActiveWindow.ScrollRow = 40
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Headless runtime may not track window scroll position; not fatal.
}
$ws.Range("G42").Select()
